# The deck ships two theme parts:
#   - theme1.xml ("Office Theme" palette) -- wired only to the Notes Master
#   - theme2.xml ("Integral" palette)     -- wired to the real Slide Master
#     (and is therefore the theme that actually paints every slide)
#
# The authored edit swaps the content of those two theme parts: the Slide
# Master's theme becomes the stock "Office Theme" colour palette, and the
# (otherwise inert) Notes-Master-only theme becomes the "Integral" palette.
#
# The PowerPoint object model only exposes a single Theme/ThemeColorScheme
# per presentation -- the one driving the Slide Master that is actually
# rendered -- so we reproduce the visible half of that swap: recolour the
# live theme's 12 scheme colours from the "Integral" palette to the
# standard "Office" palette.

$p = $ppt.ActivePresentation

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order for ThemeColorScheme.Item(n):
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb($officeThemeColors[$i - 1])
}
